$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ZQnbB960"
$ws.Range("B2").Value = 23083020
$ws.Range("C2").Value = "jdysrdy26"
$ws.Range("D2").Value = "uWXd!&97"
$ws.Range("F2").Value = "zSCmgXCC"
$ws.Range("G2").Value = "rofG"
